$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 (bold/border/centered) onto new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill I0 / IF data for rows 2-33
$data = @{
    2 = @(8, 8)
    3 = @(8, 8)
    4 = @(9, 9)
    5 = @(9, 9)
    6 = @(9, 9)
    7 = @(8, 9)
    8 = @(9, 9)
    9 = @(7, 8)
    10 = @(8, 8)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(7, 8)
    14 = @(9, 9)
    15 = @(8, 9)
    16 = @(9, 9)
    17 = @(11, 12)
    18 = @(8, 8)
    19 = @(8, 9)
    20 = @(7, 9)
    21 = @(8, 8)
    22 = @(8, 9)
    23 = @(9, 9)
    24 = @(9, 9)
    25 = @(5, 6)
    26 = @(8, 8)
    27 = @(4, 5)
    28 = @(1, 1)
    29 = @(3, 7)
    30 = @(2, 6)
    31 = @(1, 3)
    32 = @(4, 5)
    33 = @(6, 6)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
